$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 ("Naropin") already carries the exact cell-style pattern that the
# new row needs (plain/default styles for C,O,P,Q and the "center" styles
# for M,N), so clone its formatting down onto the new row 25 first and
# then overwrite the values cell by cell.
$ws.Range("A16:Q16").Copy()
$ws.Range("A25").PasteSpecial(-4122)

$ws.Range("A25").Value = 55674
$ws.Range("B25").Value = 2
$ws.Range("C25").Value = "Caverject DC 20, Injektionspräparat"
$ws.Range("D25").Value = "Pfizer AG"
$ws.Range("E25").Value = "05.99.0."
$ws.Range("F25").Value = "G04BE01"
$ws.Range("G25").Value = "Synthetika human"
$ws.Range("H25").Value = 37321
$ws.Range("I25").Value = 37321
$ws.Range("J25").Value = 42798
$ws.Range("K25").Value = 7

# Packungsgrösse ("2") is stored as text in the source file, not as a
# number - build it via a text formula on a scratch cell, copy it in as a
# value so the cell keeps its ordinary (non "quoted-text") style, then tidy
# the scratch cell back up.
$ws.Range("AA1").Formula = "=""2"""
$ws.Range("AA1").Copy()
$ws.Range("L25").PasteSpecial(-4163)
$ws.Range("AA1").Clear()

$ws.Range("M25").Value = "Spritze(n)"
$ws.Range("N25").Value = "A"
$ws.Range("O25").Value = "alprostadilum"
$ws.Range("P25").Value = "Solvens: alprostadilum 20 µg, alfadexum, lactosum anhydricum, natrii citras dihydricus, conserv.: alcohol benzylicus 4.45 mg, aqua ad iniectabilia q.s. ad solutionem pro 0.5 ml in solutione recenter reconstituta."
$ws.Range("Q25").Value = "Erektile Dysfunktion"
